$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" '28.926.49'
Set-TextValue "E2" '  +1.37%  '
Set-TextValue "D3" '1.879.12'
Set-TextValue "E3" '  -0.31%  '
Set-TextValue "D4" '1.001'
Set-TextValue "E4" '  -0.77%  '
Set-TextValue "D5" '324.95'
Set-TextValue "E5" '  -0.38%  '
Set-TextValue "D6" '1.002'
Set-TextValue "E6" '  -0.61%  '
Set-TextValue "D7" '0.4599'
Set-TextValue "E7" '  +0.36%  '
Set-TextValue "D8" '0.3876'
Set-TextValue "E8" '  +0.45%  '
Set-TextValue "D9" '0.07871'
Set-TextValue "E9" '  +0.25%  '
Set-TextValue "D10" '0.9861'
Set-TextValue "E10" '  -1.24%  '
Set-TextValue "D11" '21.76'
Set-TextValue "D12" '1.882.67'
Set-TextValue "E12" '  -0.61%  '
Set-TextValue "D13" '6.989'
Set-TextValue "E13" '  -0.90%  '
Set-TextValue "D14" '5.649'
Set-TextValue "E14" '  -0.89%  '
Set-TextValue "D15" '0.06960'
Set-TextValue "E15" '  +0.06%  '
Set-TextValue "D16" '88.00'
Set-TextValue "E16" '  +0.67%  '
Set-TextValue "E17" '  -0.81%  '
Set-TextValue "E18" '  -0.44%  '
Set-TextValue "E19" '  -0.94%  '
Set-TextValue "E20" '  -0.62%  '
Set-TextValue "D21" '28.922.86'
Set-TextValue "E21" '  +1.24%  '
Set-TextValue "E22" '  -1.39%  '
Set-TextValue "E23" '  +0.02%  '
Set-TextValue "D24" '2.101'
Set-TextValue "E24" '  +1.95%  '
Set-TextValue "D25" '156.26'
Set-TextValue "E25" '  +1.11%  '
Set-TextValue "D26" '19.30'
Set-TextValue "E26" '  -0.31%  '
Set-TextValue "D27" '6.049'
Set-TextValue "E27" '  +3.75%  '
Set-TextValue "D28" '1.928'
Set-TextValue "E28" '  -1.18%  '
Set-TextValue "D29" '117.39'
Set-TextValue "E29" '  -0.63%  '
Set-TextValue "D30" '0.09322'
Set-TextValue "E30" '  +0.07%  '
Set-TextValue "D31" '0.9025'
Set-TextValue "E31" '  -2.14%  '
Set-TextValue "D32" '5.253'
Set-TextValue "E32" '  -0.68%  '
Set-TextValue "E33" '  -0.98%  '
Set-TextValue "D34" '3.256'
Set-TextValue "E34" '  -0.32%  '
Set-TextValue "D35" '1.178'
Set-TextValue "D36" '0.05753'
Set-TextValue "E36" '  +0.03%  '
Set-TextValue "E37" '  +0.26%  '
Set-TextValue "E38" '  -0.54%  '
Set-TextValue "D39" '7.675'
Set-TextValue "E39" '  -1.34%  '
Set-TextValue "D40" '0.5648'
Set-TextValue "E40" '  -0.22%  '
Set-TextValue "E41" '  -1.27%  '
Set-TextValue "D42" '9.674'
Set-TextValue "E42" '  -0.40%  '
Set-TextValue "D43" '2.269'
Set-TextValue "E43" '  +4.43%  '
Set-TextValue "D44" '11.90'
Set-TextValue "E44" '  +1.31%  '
Set-TextValue "D45" '0.5350'
Set-TextValue "E45" '  +0.16%  '
Set-TextValue "D46" '0.07044'
Set-TextValue "E46" '  -1.43%  '
Set-TextValue "D47" '1.843'
Set-TextValue "E47" '  +0.44%  '
Set-TextValue "D48" '112.79'
Set-TextValue "E48" '  +0.30%  '
Set-TextValue "D49" '2.521'
Set-TextValue "E49" '  +2.13%  '
Set-TextValue "E50" '  -5.29%  '
Set-TextValue "D51" '70.70'
Set-TextValue "E51" '  -0.06%  '
